$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("A").Delete()
$ws.Cells.Replace("MODEL_CONDITION", "MODELCONDITION")
